# Auto-generated PowerShell-style Excel COM-interop script
# Applies the Pooh Points site update for 2026-01-27

$wb = $excel.ActiveWorkbook
$wsPlayers = $wb.Worksheets.Item("Players")
$wsOwnerTotals = $wb.Worksheets.Item("OwnerTotals")

# --- Players sheet updates ---
$wsPlayers.Range("G2").Value = "11:35 - 1st Half"
$wsPlayers.Range("H2").Value = 7
$wsPlayers.Range("I2").Value = 5
$wsPlayers.Range("K2").Value = 2
$wsPlayers.Range("O2").Value = 8
$wsPlayers.Range("G3").Value = "15:26 - 2nd Half"
$wsPlayers.Range("O3").Value = 21
$wsPlayers.Range("G4").Value = "15:26 - 2nd Half"
$wsPlayers.Range("G5").Value = "11:35 - 1st Half"
$wsPlayers.Range("J5").Value = 2
$wsPlayers.Range("O5").Value = 8
$wsPlayers.Range("G6").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H6").Value = 13
$wsPlayers.Range("I6").Value = 6
$wsPlayers.Range("J6").Value = 6
$wsPlayers.Range("N6").Value = 2
$wsPlayers.Range("O6").Value = 26
$wsPlayers.Range("G7").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H7").Value = 14
$wsPlayers.Range("O7").Value = 22
$wsPlayers.Range("G8").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H8").Value = 3
$wsPlayers.Range("I8").Value = 6
$wsPlayers.Range("O8").Value = 23
$wsPlayers.Range("G9").Value = "11:35 - 1st Half"
$wsPlayers.Range("O9").Value = 6
$wsPlayers.Range("G10").Value = "11:35 - 1st Half"
$wsPlayers.Range("O10").Value = 8
$wsPlayers.Range("G11").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H11").Value = 6
$wsPlayers.Range("L11").Value = 1
$wsPlayers.Range("O11").Value = 11
$wsPlayers.Range("G12").Value = "11:35 - 1st Half"
$wsPlayers.Range("H12").Value = -1
$wsPlayers.Range("O12").Value = 5
$wsPlayers.Range("G13").Value = "11:35 - 1st Half"
$wsPlayers.Range("G14").Value = "15:26 - 2nd Half"
$wsPlayers.Range("D15").Value = "Jadon Jones"
$wsPlayers.Range("E15").Value = "OU"
$wsPlayers.Range("F15").Value = "ARK@OU"
$wsPlayers.Range("G15").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H15").Value = 4
$wsPlayers.Range("I15").Value = 3
$wsPlayers.Range("K15").Value = 1
$wsPlayers.Range("M15").Value = 1
$wsPlayers.Range("O15").Value = 8
$wsPlayers.Range("D16").Value = "Anthony Robinson II"
$wsPlayers.Range("E16").Value = "MIZ"
$wsPlayers.Range("F16").Value = "MIZ@ALA"
$wsPlayers.Range("G16").Value = "11:35 - 1st Half"
$wsPlayers.Range("H16").Value = 1
$wsPlayers.Range("I16").Value = 3
$wsPlayers.Range("K16").Value = 0
$wsPlayers.Range("O16").Value = 4
$wsPlayers.Range("G17").Value = "11:35 - 1st Half"
$wsPlayers.Range("H17").Value = 1
$wsPlayers.Range("J17").Value = 2
$wsPlayers.Range("O17").Value = 5
$wsPlayers.Range("G18").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H18").Value = 21
$wsPlayers.Range("I18").Value = 17
$wsPlayers.Range("K18").Value = 8
$wsPlayers.Range("O18").Value = 28
$wsPlayers.Range("G19").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H19").Value = 14
$wsPlayers.Range("I19").Value = 11
$wsPlayers.Range("J19").Value = 6
$wsPlayers.Range("O19").Value = 27
$wsPlayers.Range("G20").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H20").Value = 3
$wsPlayers.Range("I20").Value = 4
$wsPlayers.Range("J20").Value = 1
$wsPlayers.Range("O20").Value = 12
$wsPlayers.Range("G21").Value = "15:26 - 2nd Half"
$wsPlayers.Range("K21").Value = 3
$wsPlayers.Range("N21").Value = 4
$wsPlayers.Range("O21").Value = 18
$wsPlayers.Range("G22").Value = "11:35 - 1st Half"
$wsPlayers.Range("O22").Value = 8
$wsPlayers.Range("G23").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H23").Value = 7
$wsPlayers.Range("I23").Value = 8
$wsPlayers.Range("O23").Value = 25
$wsPlayers.Range("G24").Value = "11:35 - 1st Half"
$wsPlayers.Range("I24").Value = 6
$wsPlayers.Range("O24").Value = 8
$wsPlayers.Range("G25").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H25").Value = 12
$wsPlayers.Range("O25").Value = 22
$wsPlayers.Range("G26").Value = "15:26 - 2nd Half"
$wsPlayers.Range("H26").Value = 10
$wsPlayers.Range("J26").Value = 2
$wsPlayers.Range("O26").Value = 14
$wsPlayers.Range("G27").Value = "11:35 - 1st Half"
$wsPlayers.Range("H27").Value = 7
$wsPlayers.Range("J27").Value = 3
$wsPlayers.Range("O27").Value = 9
$wsPlayers.Range("D28").Value = "Charles Bediako"
$wsPlayers.Range("E28").Value = "ALA"
$wsPlayers.Range("G28").Value = "11:35 - 1st Half"
$wsPlayers.Range("H28").Value = 6
$wsPlayers.Range("I28").Value = 4
$wsPlayers.Range("J28").Value = 3
$wsPlayers.Range("N28").Value = 0
$wsPlayers.Range("D29").Value = "London Jemison"
$wsPlayers.Range("G29").Value = "11:35 - 1st Half"
$wsPlayers.Range("I29").Value = 3
$wsPlayers.Range("J29").Value = 1
$wsPlayers.Range("O29").Value = 3
$wsPlayers.Range("D30").Value = "T.O. Barrett"
$wsPlayers.Range("E30").Value = "MIZ"
$wsPlayers.Range("F30").Value = "MIZ@ALA"
$wsPlayers.Range("G30").Value = "11:35 - 1st Half"
$wsPlayers.Range("H30").Value = 3
$wsPlayers.Range("I30").Value = 6
$wsPlayers.Range("J30").Value = 0
$wsPlayers.Range("D31").Value = "D.J. Wagner"
$wsPlayers.Range("E31").Value = "ARK"
$wsPlayers.Range("F31").Value = "ARK@OU"
$wsPlayers.Range("G31").Value = "15:26 - 2nd Half"
$wsPlayers.Range("J31").Value = 2
$wsPlayers.Range("N31").Value = 1
$wsPlayers.Range("O31").Value = 8
$wsPlayers.Range("G32").Value = "11:35 - 1st Half"
$wsPlayers.Range("I32").Value = 2
$wsPlayers.Range("O32").Value = 4
$wsPlayers.Range("G33").Value = "11:35 - 1st Half"

# --- OwnerTotals sheet updates ---
$wsOwnerTotals.Range("B2").Value = 21
$wsOwnerTotals.Range("B3").Value = 13
$wsOwnerTotals.Range("B4").Value = 11
$wsOwnerTotals.Range("A5").Value = "Booz"
$wsOwnerTotals.Range("B5").Value = 7
$wsOwnerTotals.Range("A6").Value = "Hal"
$wsOwnerTotals.Range("B6").Value = 4

Write-Output "Done applying PoohPoints update."
